# Update "Pais" sheet: refresh COVID country stats and fix sort order
# for Kazajistan/Moldavia, Uzbekistan/Camerun, Bulgaria/Islandia,
# Birmania/Etiopia, Belice/Nueva Caledonia and Curazao/Dominica,
# plus bump the "datos actualizados" timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1
$ws.Range("A1").Value = 'Datos actualizados a 7 de Mayo de 2020 a las 08:03'

# Row 4
$ws.Range("B4").Value = 1263197
$ws.Range("C4").Value = 105
$ws.Range("D4").Value = 213109
$ws.Range("E4").Value = 975281

# Row 17
$ws.Range("B17").Value = 53045
$ws.Range("C17").Value = 58
$ws.Range("E17").Value = 35927
$ws.Range("G17").Value = 2
$ws.Range("H17").Value = 1787

# Row 59
$ws.Range("A59").Value = 'Kazajistan'
$ws.Range("B59").Value = 4502
$ws.Range("C59").Value = 80
$ws.Range("D59").Value = 1408
$ws.Range("E59").Value = 3064
$ws.Range("F59").Value = 31
$ws.Range("H59").Value = 30

# Row 60
$ws.Range("A60").Value = 'Moldavia'
$ws.Range("B60").Value = 4476
$ws.Range("D60").Value = 1658
$ws.Range("E60").Value = 2675
$ws.Range("F60").Value = 237
$ws.Range("H60").Value = 143

# Row 72
$ws.Range("A72").Value = 'Uzbekistan'
$ws.Range("B72").Value = 2266
$ws.Range("C72").Value = 33
$ws.Range("D72").Value = 1577
$ws.Range("E72").Value = 679
$ws.Range("F72").Value = 8
$ws.Range("H72").Value = 10

# Row 73
$ws.Range("A73").Value = 'Camerun'
$ws.Range("B73").Value = 2265
$ws.Range("D73").Value = 1000
$ws.Range("E73").Value = 1157
$ws.Range("F73").Value = 12
$ws.Range("H73").Value = 108

# Row 79
$ws.Range("A79").Value = 'Bulgaria'
$ws.Range("B79").Value = 1811
$ws.Range("C79").Value = 33
$ws.Range("D79").Value = 384
$ws.Range("E79").Value = 1343
$ws.Range("F79").Value = 38
$ws.Range("H79").Value = 84

# Row 80
$ws.Range("A80").Value = 'Islandia'
$ws.Range("B80").Value = 1799
$ws.Range("D80").Value = 1750
$ws.Range("E80").Value = 39
$ws.Range("F80").Value = 0
$ws.Range("H80").Value = 10

# Row 143
$ws.Range("A143").Value = 'Birmania'
$ws.Range("C143").Value = 1
$ws.Range("D143").Value = 50
$ws.Range("E143").Value = 106
$ws.Range("H143").Value = 6

# Row 144
$ws.Range("A144").Value = 'Etiopia'
$ws.Range("B144").Value = 162
$ws.Range("D144").Value = 93
$ws.Range("E144").Value = 65
$ws.Range("H144").Value = 4

# Row 191
$ws.Range("A191").Value = 'Belice'
$ws.Range("D191").Value = 16
$ws.Range("H191").Value = 2

# Row 192
$ws.Range("A192").Value = 'Nueva Caledonia'
$ws.Range("D192").Value = 18
$ws.Range("H192").Value = 0

# Row 198
$ws.Range("A198").Value = 'Curazao'
$ws.Range("D198").Value = 13
$ws.Range("H198").Value = 1

# Row 199
$ws.Range("A199").Value = 'Dominica'
$ws.Range("D199").Value = 14
$ws.Range("H199").Value = 0
